# Updates cryptos list - applies 98 cell changes across rows 2-51
# (values from coinranking.com crypto price snapshot refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values look like plain decimal numbers (e.g. "0.989").
# Assigning such a string via .Value would make Excel auto-convert it to a
# floating point number (losing the exact text, e.g. trailing zeros / precision).
# The source column is text, so force the whole Price column to Text format
# first, assign the literal strings, then restore the cell style so the saved
# file keeps the original (default) style index on every cell.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D4').Value = '0.989'
$ws.Range('D5').Value = '211.15'
$ws.Range('D6').Value = '0.494'
$ws.Range('D8').Value = '23.08'
$ws.Range('D10').Value = '0.0597'
$ws.Range('D14').Value = '3.74'
$ws.Range('D15').Value = '0.520'
$ws.Range('D17').Value = '63.44'
$ws.Range('D18').Value = '229.36'
$ws.Range('D19').Value = '7.47'
$ws.Range('D21').Value = '0.990'
$ws.Range('D22').Value = '4.10'
$ws.Range('D23').Value = '9.29'
$ws.Range('D24').Value = '1.93'
$ws.Range('D25').Value = '150.39'
$ws.Range('D26').Value = '15.22'
$ws.Range('D28').Value = '6.56'
$ws.Range('D29').Value = '0.990'
$ws.Range('D30').Value = '1.14'
$ws.Range('D31').Value = '0.0472'
$ws.Range('D32').Value = '3.23'
$ws.Range('D33').Value = '3.12'
$ws.Range('D39').Value = '0.542'
$ws.Range('D40').Value = '2.46'
$ws.Range('D41').Value = '0.807'
$ws.Range('D42').Value = '5.65'
$ws.Range('D43').Value = '0.989'
$ws.Range('D44').Value = '1.84'
$ws.Range('D46').Value = '63.73'
$ws.Range('D48').Value = '86.44'
$ws.Range('D51').Value = '39.75'

# Restore default styling on the Price column (clears the temporary text format)
$priceRange.Style = "Normal"

# Remaining cell updates: Coin / Link text, non-ambiguous Price text, and all
# Volume(1h) percentage strings (never numeric-ambiguous because of the "%"/spaces)
$ws.Range('D2').Value = '27.827.30'
$ws.Range('E2').Value = '  +2.57%  '
$ws.Range('D3').Value = '1.568.60'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  -1.88%  '
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('E6').Value = '  +0.70%  '
$ws.Range('E7').Value = '  -2.04%  '
$ws.Range('E8').Value = '  +5.02%  '
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('D12').Value = '1.791.30'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = '1.558.95'
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').Value = '27.743.42'
$ws.Range('E16').Value = '  +2.26%  '
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('E18').Value = '  +6.86%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0703'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('E21').Value = '  -1.91%  '
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('E25').Value = '  -2.54%  '
$ws.Range('E26').Value = '  +1.20%  '
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('E29').Value = '  -1.91%  '
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').Value = '1.418.71'
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('E36').Value = '  -4.26%  '
$ws.Range('E37').Value = '  -2.27%  '
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('E39').Value = '  +2.16%  '
$ws.Range('E40').Value = '  +4.85%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E42').Value = '  -3.17%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E43').Value = '  -1.93%  '
$ws.Range('E44').Value = '  +6.03%  '
$ws.Range('E45').Value = '  -3.39%  '
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('D47').Value = '1.697.38'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('D50').Value = '0.0₇0997'
$ws.Range('E50').Value = '  -2.63%  '
$ws.Range('E51').Value = '  +17.24%  '

Write-Output "Applied 98 cell updates"
